# Apply cryptos.xlsx price/volume update (commit: "Updated cryptos list on Wed May  8 23:23:47 UTC 2024 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new values would otherwise be
# auto-detected as numbers by Excel (keeps them as text, like the originals).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"

# Row-by-row cell updates
# Row 2
$ws.Range("D2").Value = '61.010.64'
$ws.Range("E2").Value = '  -2.12%  '

# Row 3
$ws.Range("D3").Value = '2.967.14'
$ws.Range("E3").Value = '  -1.56%  '

# Row 4
$ws.Range("E4").Value = '  +0.11%  '

# Row 5
$ws.Range("D5").Value = '587.13'
$ws.Range("E5").Value = '  +1.61%  '

# Row 6
$ws.Range("D6").Value = '142.01'
$ws.Range("E6").Value = '  -4.71%  '

# Row 7
$ws.Range("E7").Value = '  +0.03%  '

# Row 8
$ws.Range("E8").Value = '  -1.33%  '

# Row 9
$ws.Range("D9").Value = '2.963.31'
$ws.Range("E9").Value = '  -1.69%  '

# Row 10
$ws.Range("E10").Value = '  -5.05%  '

# Row 11
$ws.Range("D11").Value = '5.76'
$ws.Range("E11").Value = '  +1.19%  '

# Row 12
$ws.Range("E12").Value = '  +2.49%  '

# Row 13
$ws.Range("E13").Value = '  -2.31%  '

# Row 14
$ws.Range("D14").Value = '33.88'
$ws.Range("E14").Value = '  -4.36%  '

# Row 15
$ws.Range("E15").Value = '  +1.96%  '

# Row 16
$ws.Range("D16").Value = '3.463.64'
$ws.Range("E16").Value = '  -1.42%  '

# Row 17
$ws.Range("D17").Value = '7.01'
$ws.Range("E17").Value = '  +0.06%  '

# Row 18
$ws.Range("D18").Value = '61.102.12'
$ws.Range("E18").Value = '  -1.99%  '

# Row 19
$ws.Range("D19").Value = '2.971.61'
$ws.Range("E19").Value = '  -1.39%  '

# Row 20
$ws.Range("D20").Value = '446.78'
$ws.Range("E20").Value = '  -5.21%  '

# Row 21
$ws.Range("D21").Value = '13.89'
$ws.Range("E21").Value = '  -0.98%  '

# Row 22
$ws.Range("D22").Value = '0.681'
$ws.Range("E22").Value = '  -1.70%  '

# Row 23
$ws.Range("E23").Value = '  -0.83%  '

# Row 24
$ws.Range("D24").Value = '81.31'
$ws.Range("E24").Value = '  +0.74%  '

# Row 25
$ws.Range("E25").Value = '  -2.70%  '

# Row 26
$ws.Range("D26").Value = '2.15'
$ws.Range("E26").Value = '  -8.22%  '

# Row 27
$ws.Range("D27").Value = '0.999'
$ws.Range("E27").Value = '  +0.01%  '

# Row 28
$ws.Range("D28").Value = '9.91'
$ws.Range("E28").Value = '  -5.14%  '

# Row 29
$ws.Range("E29").Value = '  +0.16%  '

# Row 30
$ws.Range("E30").Value = '  +1.18%  '

# Row 31
$ws.Range("D31").Value = '6.82'
$ws.Range("E31").Value = '  -4.96%  '

# Row 32
$ws.Range("D32").Value = '2.04'
$ws.Range("E32").Value = '  -5.64%  '

# Row 33
$ws.Range("D33").Value = '27.12'
$ws.Range("E33").Value = '  -0.17%  '

# Row 34
$ws.Range("E34").Value = '  -2.68%  '

# Row 35
$ws.Range("E35").Value = '  -3.10%  '

# Row 36
$ws.Range("D36").Value = '0.0₃0778'
$ws.Range("E36").Value = '  -1.99%  '

# Row 37
$ws.Range("E37").Value = '  -1.18%  '

# Row 38
$ws.Range("D38").Value = '9.18'
$ws.Range("E38").Value = '  +2.03%  '

# Row 39
$ws.Range("D39").Value = '50.06'
$ws.Range("E39").Value = '  -0.04%  '

# Row 40
$ws.Range("D40").Value = '2.05'
$ws.Range("E40").Value = '  -4.77%  '

# Row 41
$ws.Range("E41").Value = '  +4.98%  '

# Row 42
$ws.Range("D42").Value = '2.75'
$ws.Range("E42").Value = '  -8.76%  '

# Row 43
$ws.Range("D43").Value = '389.49'
$ws.Range("E43").Value = '  -7.04%  '

# Row 44
$ws.Range("D44").Value = '0.0351'
$ws.Range("E44").Value = '  -1.28%  '

# Row 45
$ws.Range("B45").Value = 'Maker'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D45").Value = '2.692.09'
$ws.Range("E45").Value = '  -3.92%  '

# Row 46
$ws.Range("B46").Value = 'TheGraph'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D46").Value = '0.262'
$ws.Range("E46").Value = '  -6.24%  '

# Row 47
$ws.Range("D47").Value = '37.23'
$ws.Range("E47").Value = '  -2.32%  '

# Row 48
$ws.Range("D48").Value = '131.04'
$ws.Range("E48").Value = '  +3.09%  '

# Row 49
$ws.Range("E49").Value = '  +0.08%  '

# Row 50
$ws.Range("E50").Value = '  -0.97%  '

# Row 51
$ws.Range("E51").Value = '  +0.19%  '
